$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# TreatmentTab row ("TabQuery" column, cell B5): the "Treatment Agent"
# expression used to wrap the REPLACE() call in a redundant CONCAT(...).
# Drop the CONCAT wrapper since REPLACE already returns a single scalar
# value - this was the mapping fix referenced in the commit message.
$newTreatmentQuery = @'
SELECT
    DISTINCT prt.participant_id AS "Participant Id",
    trt.treatment_id AS "Treatment Id",
    CASE 
    WHEN trt.age_at_treatment_start = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_start >= 1000 THEN 
        substr(trt.age_at_treatment_start, 1, length(trt.age_at_treatment_start) - 3) || ',' || substr(trt.age_at_treatment_start, -3)
    ELSE 
        trt.age_at_treatment_start 
END AS "Age at Treatment Start",
    CASE 
    WHEN trt.age_at_treatment_end = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_end >= 1000 THEN 
        substr(trt.age_at_treatment_end, 1, length(trt.age_at_treatment_end) - 3) || ',' || substr(trt.age_at_treatment_end, -3)
    ELSE 
        trt.age_at_treatment_end 
END AS "Age at Treatment End",
    trt.treatment_type AS "Treatment Type",
    REPLACE(trt.treatment_agent, ';', ', ') AS "Treatment Agent",
    std.dbgap_accession AS "dbGaP Accession"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnoses dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_treatments trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_reference_files rfs ON std.id = rfs."study.id"
WHERE 
    std.dbgap_accession = 'phs002599' AND prt.sex_at_birth = 'Female'
ORDER BY 
    trt.treatment_id ASC
LIMIT 100;
'@

$b5 = $ws.Range("B5")
$b5.Value2 = $newTreatmentQuery
$b5.WrapText = $true
$b5.Font.Size = 12

# The StudiesTab "StatQuery" cell (C2, the COUNT(DISTINCT ...) query) is
# left untouched - its text is unchanged by this edit.

# Reflect where the user's cursor ended up after making the edit.
$ws.Range("C5").Select() | Out-Null
